# Apply cryptos list update (Tue Oct  1 04:40:40 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "5.82", "0.1000", "63.826.17").
# Excel's COM Value setter auto-coerces such strings to numbers, which would strip
# trailing zeros / introduce float noise. Force text format for the D range first,
# assign the literal strings, then clear the (now unneeded) formatting so the cells
# end up with no style index attached -- matching their original plain state.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.826.17"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.639.47"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "581.00"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "155.49"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").Value = "2.636.24"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").Value = "5.82"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "0.383"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "28.52"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "3.111.42"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "63.792.68"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "2.637.20"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "12.15"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("D21").Value = "4.54"
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("D22").Value = "345.28"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +10.03%  "
$ws.Range("D25").Value = "67.97"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "0.0000110"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").Value = "602.95"
$ws.Range("E27").Value = "  +8.52%  "
$ws.Range("D28").Value = "9.28"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "1.61"
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("D30").Value = "8.12"
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").Value = "5.46"
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("D37").Value = "0.404"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "19.77"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "150.92"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.57"
$ws.Range("E42").Value = "  +5.50%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "41.92"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "160.55"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").Value = "24.67"
$ws.Range("E46").Value = "  +8.66%  "
$ws.Range("D47").Value = "3.92"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "0.0589"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "0.635"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "0.1000"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").Value = "0.0249"
$ws.Range("E51").Value = "  -0.31%  "

# Drop the temporary text-format style so untouched/edited D cells have no style index,
# same as before the edit (values remain text since their cell type was already fixed).
$ws.Range("D2:D51").ClearFormats()
